$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.ClearFormats()
}

# --- Row 26: now Fetch.AI (was RenderToken) ---
$ws.Cells.Item(26, 2).Value = "Fetch.AI"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell $ws.Cells.Item(26, 4) '2.23'
$ws.Cells.Item(26, 5).Value = '  +10.78%  '

# --- Row 27: now RenderToken (was Fetch.AI) ---
$ws.Cells.Item(27, 2).Value = "RenderToken"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws.Cells.Item(27, 4) '10.71'
$ws.Cells.Item(27, 5).Value = '  +6.96%  '

# --- Remaining D/E price + volume updates ---
Set-TextCell $ws.Cells.Item(2, 4) '65.445.42'
$ws.Cells.Item(2, 5).Value = '  +6.55%  '
Set-TextCell $ws.Cells.Item(3, 4) '2.993.25'
$ws.Cells.Item(3, 5).Value = '  +3.73%  '
$ws.Cells.Item(4, 5).Value = '  +0.01%  '
Set-TextCell $ws.Cells.Item(5, 4) '584.04'
Set-TextCell $ws.Cells.Item(6, 4) '153.54'
$ws.Cells.Item(6, 5).Value = '  +6.80%  '
Set-TextCell $ws.Cells.Item(7, 4) '0.999'
$ws.Cells.Item(7, 5).Value = '  -0.10%  '
Set-TextCell $ws.Cells.Item(8, 4) '0.516'
$ws.Cells.Item(8, 5).Value = '  +2.14%  '
Set-TextCell $ws.Cells.Item(9, 4) '2.989.75'
$ws.Cells.Item(9, 5).Value = '  +3.65%  '
Set-TextCell $ws.Cells.Item(10, 4) '6.99'
$ws.Cells.Item(10, 5).Value = '  -0.25%  '
Set-TextCell $ws.Cells.Item(12, 4) '0.448'
$ws.Cells.Item(12, 5).Value = '  +3.64%  '
$ws.Cells.Item(13, 5).Value = '  +3.15%  '
Set-TextCell $ws.Cells.Item(14, 4) '33.92'
$ws.Cells.Item(14, 5).Value = '  +5.96%  '
$ws.Cells.Item(15, 5).Value = '  +0.75%  '
Set-TextCell $ws.Cells.Item(16, 4) '65.309.35'
$ws.Cells.Item(16, 5).Value = '  +6.34%  '
Set-TextCell $ws.Cells.Item(17, 4) '3.488.68'
$ws.Cells.Item(17, 5).Value = '  +3.74%  '
Set-TextCell $ws.Cells.Item(18, 4) '6.92'
$ws.Cells.Item(18, 5).Value = '  +5.49%  '
Set-TextCell $ws.Cells.Item(19, 4) '2.991.09'
$ws.Cells.Item(19, 5).Value = '  +3.91%  '
Set-TextCell $ws.Cells.Item(20, 4) '452.96'
$ws.Cells.Item(20, 5).Value = '  +4.98%  '
Set-TextCell $ws.Cells.Item(21, 4) '13.72'
$ws.Cells.Item(21, 5).Value = '  +4.89%  '
$ws.Cells.Item(22, 5).Value = '  +3.77%  '
Set-TextCell $ws.Cells.Item(23, 4) '7.32'
$ws.Cells.Item(23, 5).Value = '  +7.16%  '
Set-TextCell $ws.Cells.Item(24, 4) '81.26'
$ws.Cells.Item(24, 5).Value = '  +2.50%  '
Set-TextCell $ws.Cells.Item(25, 4) '12.41'
$ws.Cells.Item(25, 5).Value = '  +3.95%  '
Set-TextCell $ws.Cells.Item(28, 4) '1.00'
$ws.Cells.Item(28, 5).Value = '  -0.07%  '
$ws.Cells.Item(29, 5).Value = '  +17.25%  '
$ws.Cells.Item(30, 5).Value = '  +11.38%  '
Set-TextCell $ws.Cells.Item(31, 4) '0.0000104'
$ws.Cells.Item(31, 5).Value = '  -0.98%  '
Set-TextCell $ws.Cells.Item(32, 4) '2.59'
$ws.Cells.Item(32, 5).Value = '  +3.77%  '
$ws.Cells.Item(33, 5).Value = '  +3.93%  '
Set-TextCell $ws.Cells.Item(34, 4) '26.84'
$ws.Cells.Item(34, 5).Value = '  +5.60%  '
Set-TextCell $ws.Cells.Item(35, 4) '0.999'
$ws.Cells.Item(35, 5).Value = '  -0.12%  '
Set-TextCell $ws.Cells.Item(36, 4) '0.984'
$ws.Cells.Item(36, 5).Value = '  +2.54%  '
Set-TextCell $ws.Cells.Item(37, 4) '5.76'
$ws.Cells.Item(37, 5).Value = '  +7.10%  '
$ws.Cells.Item(38, 5).Value = '  +9.18%  '
Set-TextCell $ws.Cells.Item(39, 4) '45.98'
$ws.Cells.Item(39, 5).Value = '  +17.10%  '
Set-TextCell $ws.Cells.Item(40, 4) '49.17'
$ws.Cells.Item(40, 5).Value = '  +0.59%  '
Set-TextCell $ws.Cells.Item(41, 4) '2.89'
$ws.Cells.Item(41, 5).Value = '  +1.97%  '
$ws.Cells.Item(42, 5).Value = '  +5.90%  '
Set-TextCell $ws.Cells.Item(43, 4) '0.297'
$ws.Cells.Item(43, 5).Value = '  +11.45%  '
$ws.Cells.Item(44, 5).Value = '  +2.43%  '
Set-TextCell $ws.Cells.Item(45, 4) '383.14'
$ws.Cells.Item(45, 5).Value = '  +11.63%  '
Set-TextCell $ws.Cells.Item(46, 4) '2.767.91'
$ws.Cells.Item(46, 5).Value = '  +2.18%  '
$ws.Cells.Item(47, 5).Value = '  +4.29%  '
Set-TextCell $ws.Cells.Item(48, 4) '134.69'
$ws.Cells.Item(48, 5).Value = '  +1.21%  '
Set-TextCell $ws.Cells.Item(50, 4) '0.106'
$ws.Cells.Item(50, 5).Value = '  +2.73%  '
Set-TextCell $ws.Cells.Item(51, 4) '23.11'
$ws.Cells.Item(51, 5).Value = '  +7.10%  '
